$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" row (row 11); everything below shifts up one row
$meta.Rows.Item(11).Delete()

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date: refreshed publication date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail" -> now "Jurisdiction" / "United States of America"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# Case Sensitive value (store as literal text "true", not boolean TRUE)
$meta.Range("B14").Value = "'true"

# --- Concepts sheet ---
# Reorder the code rows alphabetically: chat, email, mail, phone, sms
# (the "Level" column holds the text "1", not a number, so force text via the quote prefix)
$concepts = $wb.Worksheets.Item("Concepts")

$concepts.Range("A2").Value = "'1"
$concepts.Range("B2").Value = "chat"
$concepts.Range("C2").Value = "Chat"
$concepts.Range("D2").Value = "Conversational chat messaging"

$concepts.Range("A3").Value = "'1"
$concepts.Range("B3").Value = "email"
$concepts.Range("C3").Value = "Email"
$concepts.Range("D3").Value = "Email messaging"

$concepts.Range("A4").Value = "'1"
$concepts.Range("B4").Value = "mail"
$concepts.Range("C4").Value = "Mail"
$concepts.Range("D4").Value = "Postal mail messaging"

$concepts.Range("A5").Value = "'1"
$concepts.Range("B5").Value = "phone"
$concepts.Range("C5").Value = "Phone"
$concepts.Range("D5").Value = "Voice messaging by phone"

$concepts.Range("A6").Value = "'1"
$concepts.Range("B6").Value = "sms"
$concepts.Range("C6").Value = "Sms Text"
$concepts.Range("D6").Value = "SMS text messaging"
